# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice*, Leve price/profit columns) for a handful of
# leve rows across several crafting-job sheets, using freshly pulled
# market data. Two rows on the BSM sheet had stale/garbage price data
# (flat 80000) which is reset to 0, clearing their now-empty HQ profit
# cell entirely.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1746.5834
$ws.Range("I129").Value = 688.5
$ws.Range("J129").Value = 1808.8235
$ws.Range("K129").Value = 2065.5
$ws.Range("L129").Value = 5426.470499999999
$ws.Range("M129").Value = 2934.5
$ws.Range("N129").Value = -15426.4705
$ws.Range("H132").Value = 4466903
$ws.Range("I132").Value = 2533.34
$ws.Range("J132").Value = 41669984
$ws.Range("K132").Value = 7600.02
$ws.Range("L132").Value = 125009952
$ws.Range("M132").Value = -5070.02
$ws.Range("N132").Value = -125015012
$ws.Range("H137").Value = 3385.9678
$ws.Range("I137").Value = 823.5417
$ws.Range("J137").Value = 12171.429
$ws.Range("K137").Value = 2470.6251
$ws.Range("L137").Value = 36514.287
$ws.Range("M137").Value = 79.3748999999998
$ws.Range("N137").Value = -41614.287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1829.5671
$ws.Range("I61").Value = 1762.3969
$ws.Range("J61").Value = 2887.5
$ws.Range("K61").Value = 1762.3969
$ws.Range("L61").Value = 2887.5
$ws.Range("M61").Value = -1550.3969
$ws.Range("N61").Value = -3311.5
$ws.Range("H74").Value = 887.9429
$ws.Range("I74").Value = 883.75757
$ws.Range("K74").Value = 883.75757
$ws.Range("M74").Value = -9.757569999999987
$ws.Range("H77").Value = 887.9429
$ws.Range("I77").Value = 883.75757
$ws.Range("K77").Value = 4418.78785
$ws.Range("M77").Value = -50.78784999999971
$ws.Range("H122").Value = 4909.325
$ws.Range("I122").Value = 5502.7812
$ws.Range("J122").Value = 2535.5
$ws.Range("K122").Value = 16508.3436
$ws.Range("L122").Value = 7606.5
$ws.Range("M122").Value = -14058.3436
$ws.Range("N122").Value = -12506.5
$ws.Range("H132").Value = 1704.5
$ws.Range("I132").Value = 1300.625
$ws.Range("J132").Value = 2714.1875
$ws.Range("K132").Value = 3901.875
$ws.Range("L132").Value = 8142.5625
$ws.Range("M132").Value = -1371.875
$ws.Range("N132").Value = -13202.5625
$ws.Range("H136").Value = 1829.5671
$ws.Range("I136").Value = 1762.3969
$ws.Range("J136").Value = 2887.5
$ws.Range("K136").Value = 5287.1907
$ws.Range("L136").Value = 8662.5
$ws.Range("M136").Value = -2737.1907
$ws.Range("N136").Value = -13762.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H134").Value = 3217.0615
$ws.Range("I134").Value = 1978.5227
$ws.Range("K134").Value = 5935.5681
$ws.Range("M134").Value = -3400.5681

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3346.6035
$ws.Range("I31").Value = 1582.3871
$ws.Range("J31").Value = 5372.185
$ws.Range("K31").Value = 1582.3871
$ws.Range("L31").Value = 5372.185
$ws.Range("M31").Value = -1287.3871
$ws.Range("N31").Value = -5962.185
$ws.Range("H34").Value = 3346.6035
$ws.Range("I34").Value = 1582.3871
$ws.Range("J34").Value = 5372.185
$ws.Range("K34").Value = 1582.3871
$ws.Range("L34").Value = 5372.185
$ws.Range("M34").Value = -1380.3871
$ws.Range("N34").Value = -5776.185
$ws.Range("H132").Value = 1573.5593
$ws.Range("I132").Value = 1310.3334
$ws.Range("K132").Value = 3931.0002
$ws.Range("M132").Value = -1401.0002
$ws.Range("H134").Value = 554805.2
$ws.Range("I134").Value = 1105.027
$ws.Range("J134").Value = 3969289.8
$ws.Range("K134").Value = 3315.081
$ws.Range("L134").Value = 11907869.4
$ws.Range("M134").Value = -780.0810000000001
$ws.Range("N134").Value = -11912939.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 45456190
$ws.Range("I113").Value = 125001750
$ws.Range("J113").Value = 1584.2858
$ws.Range("K113").Value = 375005250
$ws.Range("L113").Value = 4752.857400000001
$ws.Range("M113").Value = -375003080
$ws.Range("N113").Value = -9092.857400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3768.2727
$ws.Range("I122").Value = 3477.8572
$ws.Range("J122").Value = 4276.5
$ws.Range("K122").Value = 10433.5716
$ws.Range("L122").Value = 12829.5
$ws.Range("M122").Value = -7983.571599999999
$ws.Range("N122").Value = -17729.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2101.125
$ws.Range("I16").Value = 2061.8125
$ws.Range("J16").Value = 2179.75
$ws.Range("K16").Value = 2061.8125
$ws.Range("L16").Value = 2179.75
$ws.Range("M16").Value = -1891.8125
$ws.Range("N16").Value = -2519.75
$ws.Range("H40").Value = 3234.7742
$ws.Range("I40").Value = 7772.875
$ws.Range("J40").Value = 1656.3043
$ws.Range("K40").Value = 7772.875
$ws.Range("L40").Value = 1656.3043
$ws.Range("M40").Value = -7636.875
$ws.Range("N40").Value = -1928.3043
$ws.Range("H68").Value = 1984.1428
$ws.Range("I68").Value = 1850.8235
$ws.Range("J68").Value = 2550.75
$ws.Range("K68").Value = 1850.8235
$ws.Range("L68").Value = 2550.75
$ws.Range("M68").Value = -1101.8235
$ws.Range("N68").Value = -4048.75
$ws.Range("H71").Value = 1984.1428
$ws.Range("I71").Value = 1850.8235
$ws.Range("J71").Value = 2550.75
$ws.Range("K71").Value = 9254.1175
$ws.Range("L71").Value = 12753.75
$ws.Range("M71").Value = -5510.1175
$ws.Range("N71").Value = -20241.75
$ws.Range("H122").Value = 8768.315000000001
$ws.Range("I122").Value = 13673.25
$ws.Range("J122").Value = 7460.3335
$ws.Range("K122").Value = 41019.75
$ws.Range("L122").Value = 22381.0005
$ws.Range("M122").Value = -38569.75
$ws.Range("N122").Value = -27281.0005
$ws.Range("H136").Value = 4252.1704
$ws.Range("I136").Value = 1829.2821
$ws.Range("J136").Value = 16063.75
$ws.Range("K136").Value = 5487.846299999999
$ws.Range("L136").Value = 48191.25
$ws.Range("M136").Value = -2937.846299999999
$ws.Range("N136").Value = -53291.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1435.4429
$ws.Range("I132").Value = 1247.6608
$ws.Range("J132").Value = 2186.5715
$ws.Range("K132").Value = 3742.9824
$ws.Range("L132").Value = 6559.7145
$ws.Range("M132").Value = -1212.9824
$ws.Range("N132").Value = -11619.7145
$ws.Range("H136").Value = 1382.0159
$ws.Range("I136").Value = 668.4909
$ws.Range("J136").Value = 6287.5
$ws.Range("K136").Value = 2005.4727
$ws.Range("L136").Value = 18862.5
$ws.Range("M136").Value = 544.5273
$ws.Range("N136").Value = -23962.5
